# Replace the data rows (2-5) with the new sensor readings, drop the old
# row 6 entirely (the sheet now only has 5 data rows), and widen four
# columns (J9/J22/J26/J27 -> spreadsheet columns 10/23/27/28) from a
# stored OOXML width of 7 to 8. Note: Excel's COM ColumnWidth property is
# offset from the raw OOXML "width" attribute by the standard ~0.83
# character padding, so ColumnWidth=7.1666... round-trips to a saved
# width of exactly 8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,34
$row2[0,0] = 45175.50694444445
$row2[0,1] = 22.58
$row2[0,2] = 15.542
$row2[0,3] = 4.221
$row2[0,4] = 47.493
$row2[0,5] = 39.284
$row2[0,6] = 17.769
$row2[0,7] = 58.8
$row2[0,8] = 27.341
$row2[0,9] = 11.61
$row2[0,10] = 17.881
$row2[0,11] = 18.828
$row2[0,12] = 19.728
$row2[0,13] = 5.673
$row2[0,14] = 17.67
$row2[0,15] = 24.849
$row2[0,16] = 14.79
$row2[0,17] = 3.779
$row2[0,18] = 2.46
$row2[0,19] = 261.617
$row2[0,20] = 49.202
$row2[0,21] = 16.31
$row2[0,22] = 32.642
$row2[0,23] = 17.025
$row2[0,24] = 2.109
$row2[0,25] = 29.266
$row2[0,26] = 14.407
$row2[0,27] = 12.944
$row2[0,28] = 15.145
$row2[0,29] = 19.485
$row2[0,30] = 3.64
$row2[0,31] = 51.902
$row2[0,32] = 9.071
$row2[0,33] = 20.391
$ws.Range("A2:AH2").Value = $row2

$row3 = New-Object 'object[,]' 1,34
$row3[0,0] = 45175.51388888889
$row3[0,1] = 12.011
$row3[0,2] = 8.336
$row3[0,3] = 1.67
$row3[0,4] = 25.482
$row3[0,5] = 21.092
$row3[0,6] = 9.452
$row3[0,7] = 39.081
$row3[0,8] = 14.543
$row3[0,9] = 6.222
$row3[0,10] = 9.433999999999999
$row3[0,11] = 10.259
$row3[0,12] = 10.653
$row3[0,13] = 3.021
$row3[0,14] = 9.398999999999999
$row3[0,15] = 13.225
$row3[0,16] = 8.106
$row3[0,17] = 1.586
$row3[0,18] = 0.931
$row3[0,19] = 135.754
$row3[0,20] = 26.391
$row3[0,21] = 8.676
$row3[0,22] = 17.422
$row3[0,23] = 9.336
$row3[0,24] = 1.074
$row3[0,25] = 18.476
$row3[0,26] = 7.663
$row3[0,27] = 6.985
$row3[0,28] = 8.167999999999999
$row3[0,29] = 10.653
$row3[0,30] = 1.294
$row3[0,31] = 35.365
$row3[0,32] = 4.776
$row3[0,33] = 10.847
$ws.Range("A3:AH3").Value = $row3

$row4 = New-Object 'object[,]' 1,34
$row4[0,0] = 45175.52083333334
$row4[0,1] = 21.139
$row4[0,2] = 15.428
$row4[0,3] = 1.482
$row4[0,4] = 45.565
$row4[0,5] = 37.69
$row4[0,6] = 16.635
$row4[0,7] = 63.137
$row4[0,8] = 25.596
$row4[0,9] = 11.285
$row4[0,10] = 16.923
$row4[0,11] = 18.367
$row4[0,12] = 19.245
$row4[0,13] = 5.313
$row4[0,14] = 16.542
$row4[0,15] = 23.485
$row4[0,16] = 13.98
$row4[0,17] = 1.123
$row4[0,18] = 0.951
$row4[0,19] = 244.473
$row4[0,20] = 46.226
$row4[0,21] = 15.269
$row4[0,22] = 30.982
$row4[0,23] = 16.493
$row4[0,24] = 2.041
$row4[0,25] = 30.734
$row4[0,26] = 13.487
$row4[0,27] = 12.022
$row4[0,28] = 14.115
$row4[0,29] = 19.218
$row4[0,30] = 0.784
$row4[0,31] = 57.003
$row4[0,32] = 8.571999999999999
$row4[0,33] = 19.09
$ws.Range("A4:AH4").Value = $row4

$row5 = New-Object 'object[,]' 1,34
$row5[0,0] = 45175.52777777778
$row5[0,1] = 0.47
$row5[0,2] = 0.06
$row5[0,3] = 0.57
$row5[0,4] = 0.73
$row5[0,5] = 0.37
$row5[0,6] = 0
$row5[0,7] = 10.52
$row5[0,8] = 0.58
$row5[0,9] = 0.35
$row5[0,10] = 0.18
$row5[0,11] = 0.36
$row5[0,12] = 0
$row5[0,13] = 0
$row5[0,14] = 0.38
$row5[0,15] = 0.72
$row5[0,16] = 0.55
$row5[0,17] = 0.7
$row5[0,18] = 0.22
$row5[0,19] = 0
$row5[0,20] = 1.78
$row5[0,21] = 0.35
$row5[0,22] = 1.27
$row5[0,23] = 0.7
$row5[0,24] = 0.05
$row5[0,25] = 4.34
$row5[0,26] = 0.31
$row5[0,27] = 0.44
$row5[0,28] = 0.45
$row5[0,29] = 0.33
$row5[0,30] = 0.5600000000000001
$row5[0,31] = 10.21
$row5[0,32] = 0.08
$row5[0,33] = 0.46
$ws.Range("A5:AH5").Value = $row5

# Row 6 no longer exists in the target sheet (A1:AH6 -> A1:AH5).
$ws.Rows.Item(6).Delete()

# Widen columns J9, J22, J26, J27 (OOXML width 7 -> 8).
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
